$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows 2-5 (years 1984-1987) and shift the remaining rows up.
$ws.Range("A2:E5").Delete(-4162)
